$wb = $excel.ActiveWorkbook

# Sheet 1: DATA_RAW - add new row 24 "FAVOK"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A24").Value = "FAVÖK"
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 0
$ws1.Range("E24").Value = -249000000
$ws1.Range("F24").Value = -77000000
$ws1.Range("G24").Value = 413000000
$ws1.Range("H24").Value = 182000000
$ws1.Range("I24").Value = 125000000
$ws1.Range("J24").Value = 128000000
$ws1.Range("K24").Value = 138000000
$ws1.Range("L24").Value = 71000000
$ws1.Range("M24").Value = 54000000
$ws1.Range("N24").Value = 105000000
$ws1.Range("O24").Value = -422000000
$ws1.Range("P24").Value = -580000000
$ws1.Range("Q24").Value = 77000000
$ws1.Range("R24").Value = -131000000
$ws1.Range("S24").Value = 135000000
$ws1.Range("T24").Value = -98000000
$ws1.Range("U24").Value = -29000000
$ws1.Range("V24").Value = 95000000
$ws1.Range("W24").Value = -330000000
$ws1.Range("X24").Value = 49000000
$ws1.Range("Y24").Value = 63000000
$ws1.Range("Z24").Value = 63000000
$ws1.Range("AA24").Value = -49000000
$ws1.Range("AB24").Value = -91000000
$ws1.Range("AC24").Value = -137000000
$ws1.Range("AD24").Value = -158000000
$ws1.Range("AE24").Value = -3000000
$ws1.Range("AF24").Value = -35000000
$ws1.Range("AG24").Value = -8000000
$ws1.Range("AH24").Value = -293000000
$ws1.Range("AI24").Value = -2000000
$ws1.Range("AJ24").Value = 45000000
$ws1.Range("AK24").Value = -1000000
$ws1.Range("AL24").Value = 119000000
$ws1.Range("AM24").Value = 28000000
$ws1.Range("AN24").Value = 164000000
$ws1.Range("AO24").Value = 153000000
$ws1.Range("AP24").Value = 150000000
$ws1.Range("AQ24").Value = 348000000
$ws1.Range("AR24").Value = 84000000
$ws1.Range("AS24").Value = 59000000
$ws1.Range("AT24").Value = 186000000
$ws1.Range("AU24").Value = 570000000
$ws1.Range("AV24").Value = 245000000
$ws1.Range("AW24").Value = 173000000
$ws1.Range("AX24").Value = 449000000
$ws1.Range("AY24").Value = 0
$ws1.Range("AZ24").Value = 662000000
$ws1.Range("BA24").Value = 831000000
$ws1.Range("BB24").Value = 948000000
$ws1.Range("BC24").Value = 0
$ws1.Range("BD24").Value = 951000000
$ws1.Range("BE24").Value = 526000000
$ws1.Range("BF24").Value = -64000000
$ws1.Range("BG24").Value = 0
$ws1.Range("BH24").Value = -145000000
$ws1.Range("BI24").Value = -20000000
$ws1.Range("BJ24").Value = 224000000
$ws1.Range("BK24").Value = 0
$ws1.Range("BL24").Value = 36000000
$ws1.Range("BM24").Value = 269000000
$ws1.Range("BN24").Value = 724000000
$ws1.Range("BO24").Value = 806000000
$ws1.Range("BP24").Value = -134000000
$ws1.Range("BQ24").Value = 1270000000

# Sheet 3: gelir tablosu (ceyreklik) - add new row 11 "FAVOK"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A11").Value = "FAVÖK"
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = 0
$ws3.Range("D11").Value = 0
$ws3.Range("E11").Value = -249000000
$ws3.Range("F11").Value = -77000000
$ws3.Range("G11").Value = 413000000
$ws3.Range("H11").Value = 182000000
$ws3.Range("I11").Value = 125000000
$ws3.Range("J11").Value = 128000000
$ws3.Range("K11").Value = 138000000
$ws3.Range("L11").Value = 71000000
$ws3.Range("M11").Value = 54000000
$ws3.Range("N11").Value = 105000000
$ws3.Range("O11").Value = -422000000
$ws3.Range("P11").Value = -580000000
$ws3.Range("Q11").Value = 77000000
$ws3.Range("R11").Value = -131000000
$ws3.Range("S11").Value = 135000000
$ws3.Range("T11").Value = -98000000
$ws3.Range("U11").Value = -29000000
$ws3.Range("V11").Value = 95000000
$ws3.Range("W11").Value = -330000000
$ws3.Range("X11").Value = 49000000
$ws3.Range("Y11").Value = 63000000
$ws3.Range("Z11").Value = 63000000
$ws3.Range("AA11").Value = -49000000
$ws3.Range("AB11").Value = -91000000
$ws3.Range("AC11").Value = -137000000
$ws3.Range("AD11").Value = -158000000
$ws3.Range("AE11").Value = -3000000
$ws3.Range("AF11").Value = -35000000
$ws3.Range("AG11").Value = -8000000
$ws3.Range("AH11").Value = -293000000
$ws3.Range("AI11").Value = -2000000
$ws3.Range("AJ11").Value = 45000000
$ws3.Range("AK11").Value = -1000000
$ws3.Range("AL11").Value = 119000000
$ws3.Range("AM11").Value = 28000000
$ws3.Range("AN11").Value = 164000000
$ws3.Range("AO11").Value = 153000000
$ws3.Range("AP11").Value = 150000000
$ws3.Range("AQ11").Value = 348000000
$ws3.Range("AR11").Value = 84000000
$ws3.Range("AS11").Value = 59000000
$ws3.Range("AT11").Value = 186000000
$ws3.Range("AU11").Value = 570000000
$ws3.Range("AV11").Value = 245000000
$ws3.Range("AW11").Value = 173000000
$ws3.Range("AX11").Value = 449000000
$ws3.Range("AY11").Value = 0
$ws3.Range("AZ11").Value = 662000000
$ws3.Range("BA11").Value = 831000000
$ws3.Range("BB11").Value = 948000000
$ws3.Range("BC11").Value = 0
$ws3.Range("BD11").Value = 951000000
$ws3.Range("BE11").Value = 526000000
$ws3.Range("BF11").Value = -64000000
$ws3.Range("BG11").Value = 0
$ws3.Range("BH11").Value = -145000000
$ws3.Range("BI11").Value = -20000000
$ws3.Range("BJ11").Value = 224000000
$ws3.Range("BK11").Value = 0
$ws3.Range("BL11").Value = 36000000
$ws3.Range("BM11").Value = 269000000
$ws3.Range("BN11").Value = 724000000
$ws3.Range("BO11").Value = 806000000
$ws3.Range("BP11").Value = -134000000
$ws3.Range("BQ11").Value = 1270000000
